# Applies the crypto price/volume update described in the commit diff.
# Rows 8 and 9 also swap their Coin/Link contents (XRP <-> LidoStakedEther).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value that must remain plain text even when it looks like a
# number (e.g. "5.62"), without leaving a residual cell style behind.
function Set-TextCell {
    param($Address, $Text)
    $rng = $ws.Range($Address)
    $rng.NumberFormat = "@"
    $rng.Value = $Text
    $rng.ClearFormats()
}

# Row 2
$ws.Range("D2").Value = "62.224.17"
$ws.Range("E2").Value = "  -2.44%  "

# Row 3
$ws.Range("D3").Value = "3.000.98"
$ws.Range("E3").Value = "  -2.51%  "

# Row 4
$ws.Range("E4").Value = "  +0.10%  "

# Row 5
Set-TextCell "D5" "581.33"
$ws.Range("E5").Value = "  -0.98%  "

# Row 6
Set-TextCell "D6" "146.66"
$ws.Range("E6").Value = "  -5.50%  "

# Row 7
$ws.Range("E7").Value = "  +0.09%  "

# Row 8
$ws.Range("B8").Value = "XRP"
$ws.Range("C8").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
Set-TextCell "D8" "0.519"
$ws.Range("E8").Value = "  -3.46%  "

# Row 9
$ws.Range("B9").Value = "LidoStakedEther"
$ws.Range("C9").Value = "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"
$ws.Range("D9").Value = "3.001.41"
$ws.Range("E9").Value = "  -2.49%  "

# Row 10
$ws.Range("E10").Value = "  -5.56%  "

# Row 11
Set-TextCell "D11" "5.62"
$ws.Range("E11").Value = "  -3.90%  "

# Row 12
Set-TextCell "D12" "0.439"
$ws.Range("E12").Value = "  -2.44%  "

# Row 13
Set-TextCell "D13" "0.0000227"
$ws.Range("E13").Value = "  -4.44%  "

# Row 14
Set-TextCell "D14" "34.59"
$ws.Range("E14").Value = "  -5.96%  "

# Row 15
$ws.Range("E15").Value = "  +1.74%  "

# Row 16
$ws.Range("D16").Value = "3.498.09"
$ws.Range("E16").Value = "  -2.38%  "

# Row 17
Set-TextCell "D17" "7.01"
$ws.Range("E17").Value = "  -2.20%  "

# Row 18
$ws.Range("D18").Value = "62.269.26"
$ws.Range("E18").Value = "  -2.16%  "

# Row 19
$ws.Range("D19").Value = "3.007.21"
$ws.Range("E19").Value = "  -2.23%  "

# Row 20
Set-TextCell "D20" "454.36"
$ws.Range("E20").Value = "  -3.60%  "

# Row 21
Set-TextCell "D21" "13.85"
$ws.Range("E21").Value = "  -3.16%  "

# Row 22
Set-TextCell "D22" "0.676"
$ws.Range("E22").Value = "  -4.13%  "

# Row 23
Set-TextCell "D23" "7.29"
$ws.Range("E23").Value = "  -2.97%  "

# Row 24
$ws.Range("E24").Value = "  -6.12%  "

# Row 25
Set-TextCell "D25" "79.78"
$ws.Range("E25").Value = "  -0.82%  "

# Row 26
Set-TextCell "D26" "12.27"
$ws.Range("E26").Value = "  -4.33%  "

# Row 27
Set-TextCell "D27" "10.09"
$ws.Range("E27").Value = "  -3.23%  "

# Row 28
Set-TextCell "D28" "0.997"
$ws.Range("E28").Value = "  -0.33%  "

# Row 29
$ws.Range("E29").Value = "  +0.22%  "

# Row 30
Set-TextCell "D30" "7.15"
$ws.Range("E30").Value = "  -2.78%  "

# Row 31
$ws.Range("E31").Value = "  -2.00%  "

# Row 32
$ws.Range("E32").Value = "  -2.18%  "

# Row 33
Set-TextCell "D33" "26.86"
$ws.Range("E33").Value = "  -0.99%  "

# Row 34
Set-TextCell "D34" "0.106"
$ws.Range("E34").Value = "  -5.16%  "

# Row 35
$ws.Range("E35").Value = "  -1.85%  "

# Row 36
$ws.Range("D36").Value = "0.0₃0789"
$ws.Range("E36").Value = "  -4.61%  "

# Row 37
Set-TextCell "D37" "5.71"
$ws.Range("E37").Value = "  -4.49%  "

# Row 38
Set-TextCell "D38" "2.11"
$ws.Range("E38").Value = "  -4.03%  "

# Row 39
Set-TextCell "D39" "50.28"
$ws.Range("E39").Value = "  -0.65%  "

# Row 40
Set-TextCell "D40" "8.95"
$ws.Range("E40").Value = "  -2.06%  "

# Row 41
Set-TextCell "D41" "2.87"
$ws.Range("E41").Value = "  -11.68%  "

# Row 42
Set-TextCell "D42" "414.92"
$ws.Range("E42").Value = "  -4.23%  "

# Row 43
$ws.Range("E43").Value = "  +0.23%  "

# Row 44
$ws.Range("E44").Value = "  -5.08%  "

# Row 45
Set-TextCell "D45" "0.0352"
$ws.Range("E45").Value = "  -2.10%  "

# Row 46
$ws.Range("D46").Value = "2.766.08"
$ws.Range("E46").Value = "  -1.52%  "

# Row 47
Set-TextCell "D47" "38.05"
$ws.Range("E47").Value = "  -4.76%  "

# Row 48
Set-TextCell "D48" "128.25"
$ws.Range("E48").Value = "  -1.31%  "

# Row 49
$ws.Range("E49").Value = "  +0.03%  "

# Row 50
Set-TextCell "D50" "0.107"
$ws.Range("E50").Value = "  -1.57%  "

# Row 51
Set-TextCell "D51" "23.64"
$ws.Range("E51").Value = "  -5.24%  "
